$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29, shifting existing rows 29-109 down to 30-110.
$ws.Rows(29).Insert()

# Populate the newly inserted row 29 with the new weekly record.
# Columns A,B,C,E,F,G,H,I,J,K,Q,T are constant across the whole sheet;
# L (Calidad) and R (Origen) keep the same values this record already had.
$ws.Range("A29").Value = 9
$ws.Range("B29").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C29").Value = "Metropolitana"
$ws.Range("D29").Value = "2022-12-15"
$ws.Range("E29").Value = 13
$ws.Range("F29").Value = "Fruta"
$ws.Range("G29").Value = 100101
$ws.Range("H29").Value = "Berries"
$ws.Range("I29").Value = 100101004
$ws.Range("J29").Value = "Frambuesa"
$ws.Range("K29").Value = "Sin especificar"
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 450
$ws.Range("N29").Value = 7000
$ws.Range("O29").Value = 7000
$ws.Range("P29").Value = 7000
$ws.Range("Q29").Value = "$/bandeja 2 kilos"
$ws.Range("R29").Value = "Provincia de Curicó"
$ws.Range("S29").Value = 3500
$ws.Range("T29").Value = 2
